$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 0.1146
$ws.Range("H5").Value = -0.2406
$ws.Range("G7").Value = 0.0006
$ws.Range("H7").Value = 0.0491
$ws.Range("I7").Value = 0.0711
$ws.Range("J7").Value = 0.0484
$ws.Range("K7").Value = 0.1798
$ws.Range("L7").Value = 0.0433
$ws.Range("M7").Value = -0.0098
$ws.Range("N7").Value = -0.0532
$ws.Range("O7").Value = -0.0965
$ws.Range("G9").Value = 0.0682
$ws.Range("H9").Value = -0.2077
$ws.Range("I9").Value = -0.1207
$ws.Range("J9").Value = -0.0848
$ws.Range("K9").Value = -0.0507
$ws.Range("L9").Value = -0.0232
$ws.Range("M9").Value = -0.0329
$ws.Range("N9").Value = -0.0644
$ws.Range("O9").Value = -0.0547
$ws.Range("G12").Value = -0.0884
$ws.Range("H12").Value = -0.7622
$ws.Range("I12").Value = -0.4149
$ws.Range("J12").Value = -0.27
$ws.Range("K12").Value = -0.1292
$ws.Range("L12").Value = -0.0661
$ws.Range("M12").Value = -0.1769
$ws.Range("G13").Value = -0.1803
$ws.Range("H13").Value = -0.1796
$ws.Range("I13").Value = -0.0664
$ws.Range("J13").Value = -0.0543
$ws.Range("K13").Value = -0.601
$ws.Range("L13").Value = -0.633
$ws.Range("M13").Value = -0.3094
$ws.Range("N13").Value = -0.1899
$ws.Range("O13").Value = -0.3258
$ws.Range("F14").Value = -0.9899
$ws.Range("G14").Value = -1.4822
$ws.Range("H14").Value = -1.019
$ws.Range("I14").Value = -0.4185
$ws.Range("J14").Value = -0.4068
$ws.Range("K14").Value = -0.3092
$ws.Range("L14").Value = -0.1444
$ws.Range("M14").Value = -0.0072
$ws.Range("N14").Value = -0.0016
$ws.Range("O14").Value = 0
$ws.Range("G15").Value = -4.1061
$ws.Range("H15").Value = -3.501
$ws.Range("I15").Value = -2.2257
$ws.Range("J15").Value = -2.1498
$ws.Range("K15").Value = -2.9398
$ws.Range("L15").Value = -1.3519
$ws.Range("M15").Value = -0.9488
$ws.Range("N15").Value = -0.4768
$ws.Range("O15").Value = -0.5114
$ws.Range("G17").Value = 0.1018
$ws.Range("H17").Value = -0.1918
$ws.Range("G21").Value = 0.0154
$ws.Range("H21").Value = 0.3123
$ws.Range("G23").Value = -0.0492
$ws.Range("H23").Value = -0.154
$ws.Range("I23").Value = -0.105
$ws.Range("J23").Value = 0.0207
$ws.Range("K23").Value = 0.0393
$ws.Range("L23").Value = 0.1934
$ws.Range("M23").Value = 0.1558
$ws.Range("N23").Value = 0.1266
$ws.Range("O23").Value = 0.0979
$ws.Range("F28").Value = -0.1281
$ws.Range("G28").Value = -0.2303
$ws.Range("H28").Value = -0.0825
$ws.Range("I28").Value = -0.0427
$ws.Range("J28").Value = -0.0331
$ws.Range("K28").Value = -0.0219
$ws.Range("L28").Value = -0.0136
$ws.Range("N28").Value = 0.0013
$ws.Range("O28").Value = 0.0014
$ws.Range("G32").Value = 0.3384
$ws.Range("H32").Value = -0.0499
$ws.Range("G34").Value = -0.0409
$ws.Range("H34").Value = 0.0632
$ws.Range("I34").Value = 0.1552
$ws.Range("J34").Value = 0.3498
$ws.Range("K34").Value = 0.4344
$ws.Range("L34").Value = 0.2602
$ws.Range("M34").Value = 0.1816
$ws.Range("N34").Value = 0.0138
$ws.Range("O34").Value = -0.0959
$ws.Range("G36").Value = 0.0699
$ws.Range("H36").Value = -0.0155
$ws.Range("I36").Value = 0.0034
$ws.Range("J36").Value = -0.0213
$ws.Range("K36").Value = 0.002
$ws.Range("L36").Value = 0.0011
$ws.Range("M36").Value = 0.0016
$ws.Range("N36").Value = 0.0029
$ws.Range("O36").Value = 0.0025
$ws.Range("G39").Value = 0.1253
$ws.Range("H39").Value = 0.3485
$ws.Range("I39").Value = 0.3561
$ws.Range("J39").Value = 0.0952
$ws.Range("K39").Value = 0.0188
$ws.Range("L39").Value = -0.175
$ws.Range("M39").Value = -0.1737
$ws.Range("G40").Value = 0.016
$ws.Range("H40").Value = 0.0025
$ws.Range("I40").Value = 0.0021
$ws.Range("J40").Value = 0.0023
$ws.Range("K40").Value = 0.021
$ws.Range("L40").Value = 0.025
$ws.Range("M40").Value = 0.0136
$ws.Range("N40").Value = 0.0089
$ws.Range("O40").Value = 0.0143
$ws.Range("F41").Value = 0.0022
$ws.Range("G41").Value = -0.0092
$ws.Range("H41").Value = 0.0099
$ws.Range("I41").Value = 0.002
$ws.Range("J41").Value = 0.0043
$ws.Range("K41").Value = 0.0058
$ws.Range("L41").Value = 0.0029
$ws.Range("M41").Value = -0.0023
$ws.Range("N41").Value = -0.0021
$ws.Range("O41").Value = -0.0002
$ws.Range("G42").Value = -0.4026
$ws.Range("H42").Value = -0.2281
$ws.Range("I42").Value = -0.1297
$ws.Range("J42").Value = -0.1827
$ws.Range("K42").Value = 0.1122
$ws.Range("L42").Value = -0.2208
$ws.Range("M42").Value = -0.323
$ws.Range("N42").Value = 1.3531
$ws.Range("O42").Value = 0.6281
$ws.Range("G44").Value = 0.229
$ws.Range("H44").Value = -0.1419
$ws.Range("G48").Value = -0.0851
$ws.Range("H48").Value = 0.2358
$ws.Range("G50").Value = -0.139
$ws.Range("H50").Value = -0.152
$ws.Range("I50").Value = -0.1135
$ws.Range("J50").Value = -0.0122
$ws.Range("K50").Value = 0.0888
$ws.Range("L50").Value = 0.2859
$ws.Range("M50").Value = 0.2279
$ws.Range("N50").Value = 0.1712
$ws.Range("O50").Value = 0.1156
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = -0.0039
$ws.Range("H55").Value = 0.0031
$ws.Range("I55").Value = 0.0023
$ws.Range("J55").Value = 0.0025
$ws.Range("K55").Value = 0.0033
$ws.Range("L55").Value = 0.0025
$ws.Range("N55").Value = 0.0012
$ws.Range("O55").Value = 0.0016
